# Generate Report for Handoff
# Swap the run's generated-file GUID (04e19470-...) for a new one
# (b5db8074-...) across the Overview / zh-cn / de-de sheets, and bump the
# handoff timestamps that were recorded for this run.

$wb = $excel.ActiveWorkbook

$oldGuid = "04e19470-3f90-4ab7-8599-98cf646f278f"
$newGuid = "b5db8074-2b82-4955-876a-2f2288b25e28"
$oldHash = "e00a306be80f4c2305c53b1c092b79eb389d2a3b"
$newHash = "509171f16d2f7e2cafde1ea33cace57f9a19070d"

# All three sheets link back to the same generated markdown file in the
# source repo; that target commit/path isn't part of this edit.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11fb0022b69077a5099b1096b31e7ee3cbc2f88f/e2e/$oldGuid.md"

# Re-point a worksheet's single hyperlink at $cellRef to the same address
# it already had, just with a refreshed display label - without leaving a
# stale duplicate <hyperlink> entry behind (TextToDisplay alone does that).
function Update-HyperlinkDisplay($ws, $cellRef, $display) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellRef), $linkAddress, "", "", $display)
    $ws.Range($cellRef).Style = "HyperLink"
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: source file name ("<guid>.md") - plain value, no hyperlink
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: path and name ("e2e\<guid>.md") - carries the hyperlink to the file
Update-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md"

# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-04 03:02:07"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: source file name, hyperlinked
Update-HyperlinkDisplay $wsZhCn "A2" "$newGuid.md"

# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-04 03:01:58"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: source file name, hyperlinked
Update-HyperlinkDisplay $wsDeDe "A2" "$newGuid.md"

# G2: Latest Handoff File
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2: Latest Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-09-04 03:02:07"

Write-Host "Handback report regenerated for $newGuid"
